$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Step 1 - Condition + Action")

# New shared strings must be interned in this exact order so they land at
# sharedStrings indices 96..100 (matching the target workbook):
#   96 context.precondition = true
#   97 Dependency: Result
#   98 Dependency: Precondition
#   99 context.precondition
#  100 context.dependency = true
$ws.Cells.Item(7, 3).Value = "context.precondition = true"
$ws.Cells.Item(8, 1).Value = "Dependency: Result"
$ws.Cells.Item(7, 1).Value = "Dependency: Precondition"
$ws.Cells.Item(8, 2).Value = "context.precondition"
$ws.Cells.Item(8, 3).Value = "context.dependency = true"
$ws.Cells.Item(7, 2).Value = $true

# Widen column A to fit the new, longer labels.
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668

# Move the active tab / selection from "Step 5 - Break" to
# "Step 1 - Condition + Action".
[void]$ws.Activate()
[void]$ws.Range("C19").Select()
